# Reassign compound property rows 2-9 to match the regenerated sample/report order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCols = @("A", "B", "C", "D")
$numCols  = @("E", "F", "G", "H", "I", "J", "K", "L", "M", "N", "O", "P", "Q", "R", "S", "T", "U", "V", "W", "X")

$rows = @(
  @{ r = 2; text = @("notvalidcomp", "unidentified", $null, $null); nums = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null) },
  @{ r = 3; text = @("dichlorobenzene", "1,4-dichlorobenzene", "C6H4Cl2", "C1=CC(=CC=C1Cl)Cl"); nums = @(147, 3.4, 6, 2, 4, 0, 0.4902448979591837, 0.4823129251700681, 0.02742857142857143, 0, 0, 6, 2, 0, 0, 0, 0.517673469387755, 0.4823129251700681, 0, 0) },
  @{ r = 4; text = @("phenol", "phenol", "C6H6O", "C1=CC=C(C=C1)O"); nums = @(94.11, 1.5, 6, 0, 6, 1, 0.765763468281798, 0, 0.06426522154925088, 0.1700031877590054, 0, 6, 0, 1, 0, 0, 0.8193178195728402, 0, 0.1807140580172139, 0) },
  @{ r = 5; text = @("palmitic acid", "hexadecanoic acid", "C16H32O2", "CCCCCCCCCCCCCCCC(=O)O"); nums = @(256.42, 6.4, 16, 0, 32, 2, 0.7494579205990172, 0, 0.125793619842446, 0.1247874580765931, 15, 0, 0, 0, 1, 0.8244793697839481, 0, 0, 0, 0.1755596287341081) },
  @{ r = 6; text = @("dodecane", "dodecane", "C12H26", "CCCCCCCCCCCC"); nums = @(170.33, 6.1, 12, 0, 26, 0, 0.846192684788352, 0, 0.1538660247754359, 0, 12, 0, 0, 0, 0, 1.000058709563788, 0, 0, 0, 0) },
  @{ r = 7; text = @("capric acid", "decanoic acid", "C10H20O2", "CCCCCCCCCC(=O)O"); nums = @(172.26, 4.1, 10, 0, 20, 2, 0.6972599558806455, 0, 0.1170323928944619, 0.1857540926506444, 9, 0, 0, 0, 1, 0.7387147335423198, 0, 0, 0, 0.2613317078834321) },
  @{ r = 8; text = @("oleic acid", "(z)-octadec-9-enoic acid", "C18H34O2", "CCCCCCCCC=CCCCCCCCC(=O)O"); nums = @(282.5, 6.5, 18, 0, 34, 2, 0.7653026548672566, 0, 0.121316814159292, 0.1132672566371681, 17, 0, 0, 0, 1, 0.8405345132743363, 0, 0, 0, 0.1593522123893805) },
  @{ r = 9; text = @("naphthalene", "naphthalene", "C10H8", "C1=CC=C2C=CC=CC2=C1"); nums = @(128.17, 3.3, 10, 0, 8, 0, 0.9371147694468284, 0, 0.06291643910431459, 0, 0, 10, 0, 0, 0, 0, 1.000031208551143, 0, 0, 0) }
)

foreach ($row in $rows) {
    $r = $row.r
    for ($i = 0; $i -lt $textCols.Length; $i++) {
        $addr = "$($textCols[$i])$r"
        $val = $row.text[$i]
        if ($null -eq $val) {
            $ws.Range($addr).ClearContents()
        } else {
            $ws.Range($addr).Value = $val
        }
    }
    for ($i = 0; $i -lt $numCols.Length; $i++) {
        $addr = "$($numCols[$i])$r"
        $val = $row.nums[$i]
        if ($null -eq $val) {
            $ws.Range($addr).ClearContents()
        } else {
            $ws.Range($addr).Value = $val
        }
    }
}

Write-Output "Rewrote compound property rows 2-9."